$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.170.02'
$ws.Range("E2").Value = '  -2.10%  '

$ws.Range("D3").Value = '2.606.95'
$ws.Range("E3").Value = '  -4.07%  '

$ws.Range("E4").Value = '  -0.07%  '

$c = $ws.Range("D5")
$c.Value = "'553.16"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.45%  '

$c = $ws.Range("D6")
$c.Value = "'154.41"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -3.02%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("E9").Value = '  -2.65%  '

$ws.Range("E10").Value = '  -3.87%  '

$c = $ws.Range("D11")
$c.Value = "'5.45"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -2.42%  '

$ws.Range("E12").Value = '  -2.22%  '

$ws.Range("D13").Value = '3.076.71'
$ws.Range("E13").Value = '  -3.94%  '

$c = $ws.Range("D14")
$c.Value = "'25.64"
$c.Style = "Normal"
$ws.Range("E14").Value = '  -3.17%  '

$ws.Range("D15").Value = '62.073.66'
$ws.Range("E15").Value = '  -2.04%  '

$ws.Range("E16").Value = '  -2.53%  '

$ws.Range("D17").Value = '2.616.51'
$ws.Range("E17").Value = '  -4.01%  '

$c = $ws.Range("D18")
$c.Value = "'11.60"
$c.Style = "Normal"
$ws.Range("E18").Value = '  -4.59%  '

$c = $ws.Range("D19")
$c.Value = "'4.53"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -3.00%  '

$c = $ws.Range("D20")
$c.Value = "'339.91"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -3.28%  '

$c = $ws.Range("D21")
$c.Value = "'6.07"
$c.Style = "Normal"
$ws.Range("E21").Value = '  -5.94%  '

$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("E23").Value = '  -2.34%  '

$c = $ws.Range("D24")
$c.Value = "'62.87"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.50%  '

$ws.Range("E25").Value = '  -0.73%  '

$c = $ws.Range("D26")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.03%  '

$c = $ws.Range("D27")
$c.Value = "'8.00"
$c.Style = "Normal"
$ws.Range("E27").Value = '  -2.46%  '

$ws.Range("D28").Value = '0.0₃0825'
$ws.Range("E28").Value = '  -7.05%  '

$c = $ws.Range("D29")
$c.Value = "'7.13"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.12%  '

$c = $ws.Range("D30")
$c.Value = "'1.34"
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.91%  '

$ws.Range("E31").Value = '  -2.90%  '

$c = $ws.Range("D32")
$c.Value = "'160.29"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -3.92%  '

$ws.Range("E33").Value = '  +0.04%  '

$c = $ws.Range("D34")
$c.Value = "'19.23"
$c.Style = "Normal"
$ws.Range("E34").Value = '  -3.11%  '

$c = $ws.Range("D35")
$c.Value = "'4.69"
$c.Style = "Normal"
$ws.Range("E35").Value = '  -3.07%  '

$c = $ws.Range("D36")
$c.Value = "'1.41"
$c.Style = "Normal"
$ws.Range("E36").Value = '  -4.74%  '

$ws.Range("E37").Value = '  -3.07%  '

$c = $ws.Range("D38")
$c.Value = "'337.23"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.88%  '

$ws.Range("E39").Value = '  -1.94%  '

$ws.Range("E40").Value = '  -6.39%  '

$c = $ws.Range("D41")
$c.Value = "'37.65"
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'3.88"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -3.52%  '

$c = $ws.Range("D43")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("D44").Value = '2.147.88'
$ws.Range("E44").Value = '  +2.00%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D45")
$c.Value = "'20.35"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -4.37%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range("D46")
$c.Value = "'0.608"
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.72%  '

$c = $ws.Range("D47")
$c.Value = "'10.97"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.83%  '

$c = $ws.Range("D48")
$c.Value = "'19.65"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -4.93%  '

$c = $ws.Range("D49")
$c.Value = "'0.0547"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -4.65%  '

$ws.Range("E50").Value = '  -2.11%  '

$ws.Range("E51").Value = '  -2.70%  '
